# Revert "Predicting PM10 hourly levels / Organizing Data Gathering and
# Exploration code" — drop the pm10_limits sheet and rename the remaining
# arima_graph sheet back to Sheet1, fixing up the chart's series formulas
# so they keep pointing at the (renamed) source sheet.

$wb = $excel.ActiveWorkbook

# Remove the pm10_limits worksheet entirely.
$wb.Worksheets.Item("pm10_limits").Delete()

# Rename the remaining sheet back to the generic default name.
$ws = $wb.Worksheets.Item("arima_graph")
$ws.Name = "Sheet1"

# The chart's series formulas still reference the old sheet name -
# repoint them at Sheet1 so the cached/linked ranges stay valid.
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$count = $chart.SeriesCollection().Count
for ($i = 1; $i -le $count; $i++) {
    $series = $chart.SeriesCollection($i)
    $series.Formula = $series.Formula -replace "arima_graph!", "Sheet1!"
}

$ws.Activate()
